$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text values are preserved exactly (no auto numeric/date conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.434.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3838"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07833"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9873"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.83"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.897"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.634"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06966"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009974"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.467.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.256"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.190.73"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.43%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.107"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.639"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.939"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09265"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9072"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.264"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.317"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.302"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.140"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02056"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.690"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5561"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1770"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.639"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07102"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5218"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.139"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.120"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.811"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.409"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.25%  "
